# Applies the betexplorer scraper refresh for croatia/prva-nl 2023-2024:
#  - several existing rows had their match data (columns F:V) reshuffled
#    between rows sharing the same match-day (A:E index/date columns stay put)
#  - four brand-new match rows were appended at the bottom (75-78)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowFV($row) {
    $vals = @()
    for ($c = 6; $c -le 22; $c++) {
        $vals += ,($ws.Cells.Item($row, $c).Value())
    }
    return $vals
}

function Set-RowFV($row, $vals) {
    for ($i = 0; $i -lt 17; $i++) {
        $c = 6 + $i
        $ws.Cells.Item($row, $c).Value = $vals[$i]
    }
}

# --- Reshuffle match data (columns F:V only) among existing rows ---------

# Rows 3,4,5,6: new3=old4, new4=old5, new5=old6, new6=old3
$r3 = Get-RowFV 3
$r4 = Get-RowFV 4
$r5 = Get-RowFV 5
$r6 = Get-RowFV 6
Set-RowFV 3 $r4
Set-RowFV 4 $r5
Set-RowFV 5 $r6
Set-RowFV 6 $r3

# Rows 40,41: swap
$r40 = Get-RowFV 40
$r41 = Get-RowFV 41
Set-RowFV 40 $r41
Set-RowFV 41 $r40

# Rows 52,53: swap
$r52 = Get-RowFV 52
$r53 = Get-RowFV 53
Set-RowFV 52 $r53
Set-RowFV 53 $r52

# Rows 57,58,59: new57=old59, new58=old57, new59=old58
$r57 = Get-RowFV 57
$r58 = Get-RowFV 58
$r59 = Get-RowFV 59
Set-RowFV 57 $r59
Set-RowFV 58 $r57
Set-RowFV 59 $r58

# Rows 70,71: swap
$r70 = Get-RowFV 70
$r71 = Get-RowFV 71
Set-RowFV 70 $r71
Set-RowFV 71 $r70

# --- Append four new match rows (75-78) -----------------------------------

$ws.Range("A2").Copy()
$ws.Range("A75:A78").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E75:E78").PasteSpecial(-4122)

$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = "croatia"
$ws.Cells.Item(75, 3).Value = "prva-nl"
$ws.Cells.Item(75, 4).Value = "2023-2024"
$ws.Cells.Item(75, 5).Value = 45234.57638888889
$ws.Cells.Item(75, 6).Value = "Zrinski Jurjevac"
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = "Vukovar 1991"
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 2.12
$ws.Cells.Item(75, 11).Value = "03/11/2023 02:12"
$ws.Cells.Item(75, 12).Value = 2.06
$ws.Cells.Item(75, 13).Value = "04/11/2023 13:40"
$ws.Cells.Item(75, 14).Value = 3.23
$ws.Cells.Item(75, 15).Value = "03/11/2023 02:12"
$ws.Cells.Item(75, 16).Value = 3.33
$ws.Cells.Item(75, 17).Value = "04/11/2023 13:40"
$ws.Cells.Item(75, 18).Value = 3.01
$ws.Cells.Item(75, 19).Value = "03/11/2023 02:12"
$ws.Cells.Item(75, 20).Value = 3.5
$ws.Cells.Item(75, 21).Value = "04/11/2023 13:40"
$ws.Cells.Item(75, 22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/zrinski-jurjevac-vukovar-1991/vwSNvvA2/"

$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = "croatia"
$ws.Cells.Item(76, 3).Value = "prva-nl"
$ws.Cells.Item(76, 4).Value = "2023-2024"
$ws.Cells.Item(76, 5).Value = 45234.58333333334
$ws.Cells.Item(76, 6).Value = "Cibalia"
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = "Dubrava"
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 1.98
$ws.Cells.Item(76, 11).Value = "03/11/2023 02:12"
$ws.Cells.Item(76, 12).Value = 2.56
$ws.Cells.Item(76, 13).Value = "04/11/2023 13:59"
$ws.Cells.Item(76, 14).Value = 3.32
$ws.Cells.Item(76, 15).Value = "03/11/2023 02:12"
$ws.Cells.Item(76, 16).Value = 3.28
$ws.Cells.Item(76, 17).Value = "04/11/2023 13:59"
$ws.Cells.Item(76, 18).Value = 3.25
$ws.Cells.Item(76, 19).Value = "03/11/2023 02:12"
$ws.Cells.Item(76, 20).Value = 2.67
$ws.Cells.Item(76, 21).Value = "04/11/2023 13:59"
$ws.Cells.Item(76, 22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/cibalia-dubrava-zagreb/nquxHzIe/"

$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = "croatia"
$ws.Cells.Item(77, 3).Value = "prva-nl"
$ws.Cells.Item(77, 4).Value = "2023-2024"
$ws.Cells.Item(77, 5).Value = 45234.58333333334
$ws.Cells.Item(77, 6).Value = "Dugopolje"
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = "Bijelo Brdo"
$ws.Cells.Item(77, 9).Value = 1
$ws.Cells.Item(77, 10).Value = 1.62
$ws.Cells.Item(77, 11).Value = "03/11/2023 02:12"
$ws.Cells.Item(77, 12).Value = 1.92
$ws.Cells.Item(77, 13).Value = "04/11/2023 13:59"
$ws.Cells.Item(77, 14).Value = 3.57
$ws.Cells.Item(77, 15).Value = "03/11/2023 02:12"
$ws.Cells.Item(77, 16).Value = 3.03
$ws.Cells.Item(77, 17).Value = "04/11/2023 13:59"
$ws.Cells.Item(77, 18).Value = 4.6
$ws.Cells.Item(77, 19).Value = "03/11/2023 02:12"
$ws.Cells.Item(77, 20).Value = 4.54
$ws.Cells.Item(77, 21).Value = "04/11/2023 13:59"
$ws.Cells.Item(77, 22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-bijelo-brdo/jHoTIdmq/"

$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = "croatia"
$ws.Cells.Item(78, 3).Value = "prva-nl"
$ws.Cells.Item(78, 4).Value = "2023-2024"
$ws.Cells.Item(78, 5).Value = 45234.58333333334
$ws.Cells.Item(78, 6).Value = "Orijent"
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = "Jarun"
$ws.Cells.Item(78, 9).Value = 1
$ws.Cells.Item(78, 10).Value = 1.79
$ws.Cells.Item(78, 11).Value = "03/11/2023 02:12"
$ws.Cells.Item(78, 12).Value = 1.79
$ws.Cells.Item(78, 13).Value = "04/11/2023 13:52"
$ws.Cells.Item(78, 14).Value = 3.62
$ws.Cells.Item(78, 15).Value = "03/11/2023 02:12"
$ws.Cells.Item(78, 16).Value = 4.01
$ws.Cells.Item(78, 17).Value = "04/11/2023 13:52"
$ws.Cells.Item(78, 18).Value = 3.55
$ws.Cells.Item(78, 19).Value = "03/11/2023 02:12"
$ws.Cells.Item(78, 20).Value = 3.79
$ws.Cells.Item(78, 21).Value = "04/11/2023 13:52"
$ws.Cells.Item(78, 22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/orijent-jarun/QZwtGfX1/"
